$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.284177422523499
$ws.Range("B1").Value = 1.56017279624939
$ws.Range("C1").Value = 1.94821572303772
$ws.Range("D1").Value = 1.781768202781677
$ws.Range("E1").Value = 1.475476145744324
